$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New event rows appended to the event list (rows 55-61)
$dates = @(46004, 46005, 46053, 46054, 46055, 46056, 46057)
$names = @("アイマス", "アイマス", "なにわ男子", "なにわ男子", "なにわ男子", "なにわ男子", "なにわ男子")

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 55 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = "🔴"
    $ws.Cells.Item($row, 2).Font.Name = "Segoe UI Symbol"
    $ws.Cells.Item($row, 3).Value = $names[$i]
}

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection / scroll position as last left by the author
$excel.ActiveWindow.ScrollRow = 48
[void]$ws.Range("C53").Select()
